$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.759.39'
$ws.Range("E2").Value = '  -5.66%  '

$ws.Range("D3").Value = '1.809.93'
$ws.Range("E3").Value = '  -4.94%  '

$ws.Range("D4").Value = '''0.9997'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '''275.24'
$ws.Range("E5").Value = '  -10.20%  '

$ws.Range("D6").Value = '''0.9996'
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Value = '''0.5073'
$ws.Range("E7").Value = '  -6.21%  '

$ws.Range("D8").Value = '''0.3517'
$ws.Range("E8").Value = '  -7.70%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '''0.06649'
$ws.Range("E9").Value = '  -8.93%  '

$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").Value = '''20.03'
$ws.Range("E10").Value = '  -9.33%  '

$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = '''0.8338'
$ws.Range("E11").Value = '  -7.53%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '''0.07807'
$ws.Range("E12").Value = '  -4.71%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.782.00'
$ws.Range("E13").Value = '  -7.41%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''5.064'
$ws.Range("E14").Value = '  -5.64%  '

$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '''87.51'
$ws.Range("E15").Value = '  -8.41%  '

$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '''0.9991'
$ws.Range("E16").Value = '  -0.14%  '

$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = '''13.87'
$ws.Range("E17").Value = '  -6.89%  '

$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").Value = '''0.9998'
$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '''0.000008000'
$ws.Range("E19").Value = '  -7.62%  '

$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("C20").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D20").Value = '25.817.05'
$ws.Range("E20").Value = '  -5.52%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''4.728'
$ws.Range("E21").Value = '  -6.39%  '

$ws.Range("B22").Value = 'Cosmos'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D22").Value = '''9.979'
$ws.Range("E22").Value = '  -8.02%  '

$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = '''6.045'
$ws.Range("E23").Value = '  -7.20%  '

$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").Value = '''141.68'
$ws.Range("E24").Value = '  -4.78%  '

$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D25").Value = '''2.197'
$ws.Range("E25").Value = '  -4.97%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '''1.654'
$ws.Range("E26").Value = '  -5.90%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''16.99'
$ws.Range("E27").Value = '  -7.61%  '

$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").Value = '''108.63'
$ws.Range("E28").Value = '  -6.82%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '''4.330'
$ws.Range("E29").Value = '  -10.67%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '''4.209'
$ws.Range("E30").Value = '  -9.81%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '''0.08779'
$ws.Range("E31").Value = '  -4.38%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.04849'
$ws.Range("E32").Value = '  -4.31%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''0.7329'
$ws.Range("E33").Value = '  -11.36%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '''2.905'
$ws.Range("E34").Value = '  -3.85%  '

$ws.Range("D35").Value = '''1.135'
$ws.Range("E35").Value = '  -7.48%  '

$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = '''0.9989'
$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = '''3.033'
$ws.Range("E37").Value = '  -8.58%  '

$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").Value = '''0.5220'
$ws.Range("E38").Value = '  -13.00%  '

$ws.Range("D39").Value = '''0.01854'
$ws.Range("E39").Value = '  -7.24%  '

$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '''2.300'
$ws.Range("E40").Value = '  -14.27%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '''0.9502'
$ws.Range("E41").Value = '  -11.67%  '

$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '''112.47'
$ws.Range("E42").Value = '  -2.95%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''6.180'
$ws.Range("E43").Value = '  -7.25%  '

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '''8.067'
$ws.Range("E44").Value = '  -13.25%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '''0.9992'
$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.4555'
$ws.Range("E46").Value = '  -11.11%  '

$ws.Range("D47").Value = '''0.1376'
$ws.Range("E47").Value = '  -10.31%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''9.279'
$ws.Range("E48").Value = '  -8.94%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''36.14'
$ws.Range("E49").Value = '  -5.25%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '''1.496'
$ws.Range("E50").Value = '  -8.56%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.05789'
$ws.Range("E51").Value = '  -5.00%  '
